# Rotate the "Onshore wind plants" / "Offshore wind plants" / "Photovoltaic plants"
# rows (rows 5/6/7 of every year sheet) by one position:
#   new row 5 <- old row 7
#   new row 6 <- old row 5
#   new row 7 <- old row 6
# This matches the shared-string reorder (Photovoltaic moved ahead of Onshore/Offshore)
# together with the corresponding rotation of the Column-E allocation values seen
# across every year worksheet (2000-2100).

$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
  $ws = $wb.Worksheets.Item($i)

  $c5 = $ws.Range("C5").Value()
  $c6 = $ws.Range("C6").Value()
  $c7 = $ws.Range("C7").Value()

  $e5 = $ws.Range("E5").Value()
  $e6 = $ws.Range("E6").Value()
  $e7 = $ws.Range("E7").Value()

  $ws.Range("C5").Value = $c7
  $ws.Range("C6").Value = $c5
  $ws.Range("C7").Value = $c6

  $ws.Range("E5").Value = $e7
  $ws.Range("E6").Value = $e5
  $ws.Range("E7").Value = $e6
}
